$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: average of the |S*|/n column (J) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 12

# --- Rows 14-17: summary labels + stats ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# bold + vertically-centered labels/values for the summary block
$ws.Range("B14:B17").Font.Bold = $true
$ws.Range("B14:B17").VerticalAlignment = -4108

# selection ends up on J12, matching the source workbook
$ws.Range("J12").Select()

# print setup (A4 portrait), matching the resaved workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
